# The document consists of a single paragraph whose only run holds an
# inline picture (a <w:drawing> wrapping a <pic:pic>). The edit removes
# that picture, leaving the paragraph empty (the paragraph mark itself
# stays, so the document keeps its one empty paragraph + sectPr).
$d = $word.ActiveDocument

if ($d.InlineShapes.Count -gt 0) {
    # Delete from the end backwards in case there were several.
    for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
        $d.InlineShapes($i).Delete()
    }
}
